$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 20, shifting existing rows 20-72 down to 21-73.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new data record.
$ws.Range("A20").Value = 7
$ws.Range("B20").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C20").Value = "Ñuble"
$ws.Range("D20").Value = 44575
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 100112030
$ws.Range("G20").Value = "Poroto granado"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = 24000
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = 24500
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Provincia de Diguillín"
$ws.Range("P20").Value = 980
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"

# Match the date style used by the rest of column D.
$ws.Range("D20").NumberFormat = $ws.Range("D21").NumberFormat
